$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45981
$ws.Range("B2").Value = 83.53
$ws.Range("C2").Value = 75.75
$ws.Range("D2").Value = 74.18000000000001
$ws.Range("E2").Value = 72.42
$ws.Range("F2").Value = 72.87
$ws.Range("G2").Value = 83.54000000000001
$ws.Range("H2").Value = 94.7
$ws.Range("I2").Value = 106.43
$ws.Range("J2").Value = 116.81
$ws.Range("K2").Value = 91.05
$ws.Range("L2").Value = 70.09
$ws.Range("M2").Value = 62.71
$ws.Range("N2").Value = 56.15
$ws.Range("O2").Value = 46.35
$ws.Range("P2").Value = 34.88
$ws.Range("Q2").Value = 45.11
$ws.Range("R2").Value = 77.19
$ws.Range("S2").Value = 102.88
$ws.Range("T2").Value = 116.58
$ws.Range("U2").Value = 124.39
$ws.Range("V2").Value = 125.1
$ws.Range("W2").Value = 105.72
$ws.Range("X2").Value = 100.96
$ws.Range("Y2").Value = 89.25
$ws.Range("Z2").Value = 84.53
$ws.Range("AB2").Value = 105.26
$ws.Range("AD2").Value = 120.48
$ws.Range("AF2").Value = 115.41
$ws.Range("AG2").Value = "0h-16h"
